# edit.ps1 - applies the "SAR/CAR/IAR" syllabus revision described in the
# commit "Added the notion of SAR, CAR, IAR to the syllabus"
#
# Strategy:
#  1. Replace the text of the three "Learning Objectives" bullet paragraphs
#     (Application Context / Compositional Rules / Technical Vocabulary) with
#     their revised wording (Value Identification.../Technological
#     Vocabulary.../Application Context...) including the new SAR/CAR/IAR
#     call-outs.
#  2. Insert two blank ListParagraph bullet paragraphs and a new summary
#     paragraph describing the SAR/CAR/IAR reports, right after the third
#     bullet.
#  3. Move the "_GoBack" bookmark from the "Content Detail" heading to the
#     end of the new summary paragraph (around the word "three").

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Paragraph 1: "Application Context: Know the broad history..." ->
#    "Value Identification and Standards formulation: Understand..."
# ---------------------------------------------------------------------
$old1 = "Application Context: Know the broad history, trend, and major industry standards of computation through a narrative integrated with hands-on projects or in-person conversation with original inventors, scholars, domain experts, and policy makers. "
$new1 = "Value Identification and Standards formulation: Understand how to frame the question of value proposition in the historical context of computing industry. Learn about the historical trends, major industry standards, and meet up and discuss the industry trends with critical personalities in person and in their direct descendants. Students are expected to create a Strategy Analysis Report (SAR) at the end of the study."

$found = $d.Content.Find.Execute($old1, $true, $true, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found) {
    throw "Paragraph 1 (Application Context) not found"
}

# ---------------------------------------------------------------------
# 2. Paragraph 2: "Compositional Rules: Understand the recombinatorial..."
#    -> "Technological Vocabulary and Compositional Rules: Learn to use
#    MediaWiki... combinatorial possibilities of the known spectrum of..."
# ---------------------------------------------------------------------
$old2 = "Compositional Rules: Understand the recombinatorial possibilities of well-known computational models, software tools, hardware technologies, and service providers of computation. How different symbolic systems can be put together to analyze, predict or control certain systems.  At the same time, we will learn about the basic notion of the Correctness of System Design, Design by Contract, Algebra of Computer Programs, the Composition of Distributed and Centralized Computing, and present computation results using Human-Machine Interface Technologies. "
$new2 = "Technological Vocabulary and Compositional Rules: Learn to use MediaWiki, Git, GitHub, Wolfram|Alpha, and other open-sourced computational services and tools for vocabulary management. Understand the combinatorial possibilities of the known spectrum of computational models, software tools, hardware technologies, and service providers of computation. How different symbolic systems can be put together to analyze, predict or control certain systems.  At the same time, we will learn about the basic notion of the Correctness of System Design, Design by Contract, Algebra of Computer Programs, the Composition of Distributed and Centralized Computing, and present computation results using Human-Machine Interface Technologies. Students are expected to produce a Computational Resource Analysis Report (CAR) at the end of study."

$found = $d.Content.Find.Execute($old2, $true, $true, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found) {
    throw "Paragraph 2 (Compositional Rules) not found"
}

# ---------------------------------------------------------------------
# 3. Paragraph 3: "Technical Vocabulary: Learn to use  MediaWiki..." ->
#    "Application Context: Given the personal interests..."
# ---------------------------------------------------------------------
$old3 = "Technical Vocabulary: Learn to use  MediaWiki, Git, GitHub/GItLab, and other collective knowledge management tools to organize the vocabulary of computational practices. "
$new3 = "Application Context: Given the personal interests and group decisions, students will identify up-to-date applications of computational thinking to their selected subjects. These applications will follow the reasoning framework of computational thinking and documented using computable languages. Students are expected to product an Industry Analysis Report (IAR) at the end of study."

$found = $d.Content.Find.Execute($old3, $true, $true, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $found) {
    throw "Paragraph 3 (Technical Vocabulary) not found"
}

Write-Output "Step 1-3 complete: bullet paragraphs updated"

# ---------------------------------------------------------------------
# 4. Insert two blank ListParagraph bullets (ind left=1680, no numbering)
#    plus the new summary paragraph (ind left=720) right after the third
#    bullet ("...Industry Analysis Report (IAR) at the end of study.").
# ---------------------------------------------------------------------
# Locate the third-bullet paragraph by its (now updated) text, via the
# Paragraphs collection so we can reliably chain inserts by index.
$thirdIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Application Context: Given the personal interests*") {
        $thirdIndex = $i
        break
    }
}
if ($thirdIndex -eq -1) {
    throw "Could not locate third bullet paragraph for insertion"
}

$thirdPara = $d.Paragraphs.Item($thirdIndex)
$thirdPara.Range.InsertParagraphAfter() | Out-Null

$blank1 = $d.Paragraphs.Item($thirdIndex + 1)
$blank1.Range.ListFormat.RemoveNumbers()
$blank1.Range.ParagraphFormat.LeftIndent = 84
$blank1.Range.InsertParagraphAfter() | Out-Null

$blank2 = $d.Paragraphs.Item($thirdIndex + 2)
$blank2.Range.ListFormat.RemoveNumbers()
$blank2.Range.ParagraphFormat.LeftIndent = 84
$blank2.Range.InsertParagraphAfter() | Out-Null

$summaryPara = $d.Paragraphs.Item($thirdIndex + 3)
$summaryPara.Range.ListFormat.RemoveNumbers()
$summaryPara.Range.ParagraphFormat.Style = "Normal"
$summaryPara.Range.ParagraphFormat.LeftIndent = 36

$apost = [char]0x2019
$summaryText = "The learning outcomes should be integrated using a digital publishing process. Students will be collecting their ideas using a MediaWiki-like workflow, to capture their ideas, and publish their thought processes after each learning session. Then, the three main threads of this study will produce three respective reports (SAR, CAR, and IAR). Each project is considered to have an incremental contribution, after a complete report is being edited, refined, and authorized to publish. SAR should contain testable statements or test cases to help verify and validate computational models or industry-specific applications. CAR is a report that articulate how to perform the computation using existing tools and services. IAR is the specific report explaining how an industry would utilize computational thinking and what are the known and expected results. Students$($apost) learning outcomes will be judged based on their contribution and the quality of the three reports."

$summaryPara.Range.InsertAfter($summaryText)

Write-Output "Step 4 complete: new paragraphs inserted"

# ---------------------------------------------------------------------
# 5. Move the "_GoBack" bookmark from the "Content Detail" heading to
#    the end of the new summary paragraph, collapsed right between
#    "...quality of the three" and " reports.".
# ---------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$matchRange = $d.Content
$found = $matchRange.Find.Execute("quality of the three reports.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate 'quality of the three reports.' text for bookmark placement"
}
$bmStart = $matchRange.Start + ("quality of the three").Length
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Step 5 complete: _GoBack bookmark relocated"
